# Aggiunta ultima tabella al dominio
#
# Appends the closing "Risorse-argomenti" entity description (one intro
# paragraph + two bulleted field rows) right after the existing trailing
# blank paragraph and before the section break, mirroring the pattern
# already used for the other entities (e.g. "Tutor-argomenti") earlier in
# the document.

$d = $word.ActiveDocument

# The document currently ends with a lone empty paragraph immediately
# before the section break. Create a fresh paragraph after it; this is
# where the new content will be injected so the existing empty paragraph
# is left untouched, exactly as in the target revision.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($n)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$target = $newPara.Range
$target.Collapse(1)

# Insert the three new paragraphs (intro sentence + two PK/FK field rows)
# as literal WordprocessingML so the run layout, bold run, list formatting
# (reusing numId 3, same as the other entities) and the proofing marks
# (w:proofErr) exactly match what Word itself produces when typing this
# text with spelling/grammar checking on.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:r><w:t xml:space="preserve">L&#8217;entit&#224; </w:t></w:r>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Risorse-argomenti </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>deve contenere</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> i seguenti dati:</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragrafoelenco"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>id_risorse</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>int</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> (5) PK, FK</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragrafoelenco"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>id_argomenti</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>int</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> (5) PK, FK</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$target.InsertXML($xml)

# InsertXML populated the scaffold paragraph created above and pushed a new
# empty paragraph after it (the OOXML rule that there must always be a
# paragraph mark after a block insertion at the very end of the body).
# Merge that leftover empty paragraph away so the new content flows
# straight into the section break, matching the target structure exactly.
$trailing = $d.Paragraphs($d.Paragraphs.Count)
$cleanup = $d.Range($trailing.Range.Start - 1, $trailing.Range.End)
$cleanup.Delete()
